$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 3" shape (Presented By / name / college block) by name
# so this keeps working even if shape ordering ever shifts.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 3") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(3)
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 2 holds the presenter's name: "1.AVANTHIKA.M" -> "1.DHARANI S"
$para2 = $tr.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "1.DHARANI S"

# Paragraph 3 holds the college name, originally split across three runs
# (" VIVEKANANDHA COLLEGE OF TECHNOLOGY " / "FOR WOMEN-CIVIL " / "ENGINEERING").
# Merge it into a single run carrying the combined text, reusing the first
# run's formatting. Route through an intermediate value so the engine treats
# it as a genuine edit (identical-text reassignment is a no-op) and
# collapses the run list.
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "ZZZ"
$para3.Text = " VIVEKANANDHA COLLEGE OF TECHNOLOGY FOR WOMEN-CIVIL ENGINEERING"
